# adding new progress as of date 04 nov 2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Delete the two obsolete SOP rows (old row 12: "Equipment Operation
# Procedure (SOP-031)" and old row 13: "Equipment  Request &handover
# procedure(SOP-028)"). Deleting shifts rows 14-17 up to 12-15, carrying
# each row's own formatting (fill/border) along with it.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()

# Narrow the TRAININGS column by one character.
$ws.Columns.Item(2).ColumnWidth = 55

# Renumber the SN column (A) for the rows that shifted up so it stays
# sequential (10, 11, 12, 13) instead of keeping the old 12, 13, 14, 15.
$sn = @{
    12 = 10
    13 = 11
    14 = 12
    15 = 13
}
foreach ($row in $sn.Keys) {
    $ws.Cells.Item($row, 1).Value = $sn[$row]
}

# Refresh "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I) for every data row
# now that the report was regenerated a day later (03-Nov-2025 ->
# 04-Nov-2025): every period-to-expire count drops by 1 day.
$periods = @{
    3 = 687
    4 = 255
    5 = 288
    6 = 262
    7 = 262
    8 = 265
    9 = 265
    10 = 297
    11 = 36
    12 = -97
    13 = 128
    14 = 140
    15 = 302
}

foreach ($row in $periods.Keys) {
    $ws.Cells.Item($row, 8).Value = $periods[$row]
    # Pre-format as text so the DD-MMM-YYYY-looking string is kept
    # literally instead of being auto-converted to a date serial.
    $ws.Cells.Item($row, 9).NumberFormat = "@"
    $ws.Cells.Item($row, 9).Value = "04-Nov-2025"
}
